$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 content: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Leave selection on the edited cell, matching the saved view state
$ws.Range("E8").Select()
